$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PERMISOS")

# Row 7 was a blank spacer row. Remove it (rows 8-14 shift up to 7-13).
$ws.Rows("7").Delete()

# Make room for the new "CLIENTES / GESTION" entry (row 8) and a blank gap (row 9),
# pushing the MANTENIMIENTO block back down so it starts at row 10.
$ws.Rows("7:9").Insert()

# Row 7 and row 9 stay empty (fully cleared, including the inherited style).
$ws.Range("C7").Clear()
$ws.Range("C9").Clear()

# Row 8: new CLIENTES / GESTION permission entry.
$ws.Range("B8").Value = "CLIENTES"
$ws.Range("C8").Value = "GESTION"
$ws.Range("E8").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B8,""'"","","",""'"",C8,""'"","")"","";"")"

# Fix the table name used by every remaining generated-SQL formula:
# "permisos (modulo,submodulo)" -> "permiso(modulo,submodulo)"
$ws.Range("E3").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B3,""'"","","",""'"",C3,""'"","")"","";"")"
$ws.Range("E5").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B5,""'"","","",""'"",C5,""'"","")"","";"")"
$ws.Range("E6").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B6,""'"","","",""'"",C6,""'"","")"","";"")"
$ws.Range("E10").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B10,""'"","","",""'"",C10,""'"","")"","";"")"
$ws.Range("E11").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B11,""'"","","",""'"",C11,""'"","")"","";"")"
$ws.Range("E12").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B12,""'"","","",""'"",C12,""'"","")"","";"")"
$ws.Range("E13").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B13,""'"","","",""'"",C13,""'"","")"","";"")"
$ws.Range("E14").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B14,""'"","","",""'"",C14,""'"","")"","";"")"

# Two new rows continuing the MANTENIMIENTO block.
$ws.Range("B15").Value = "MANTENIMIENTO"
$ws.Range("C15").Value = "CREDITOS_TIPOS"
$ws.Range("C15").Interior.Color = $ws.Range("C14").Interior.Color
$ws.Range("E15").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B15,""'"","","",""'"",C15,""'"","")"","";"")"

$ws.Range("B16").Value = "MANTENIMIENTO"
$ws.Range("C16").Value = "CREDITOS_ESTADOS"
$ws.Range("C16").Interior.Color = $ws.Range("C14").Interior.Color
$ws.Range("E16").Formula = "=CONCATENATE(""insert into permiso(modulo,submodulo) values ("",""'"",B16,""'"","","",""'"",C16,""'"","")"","";"")"

$ws.Range("E18").Select()
